$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "0.9993") must be
# force-formatted as Text first, otherwise Excel auto-converts the literal
# into a real number, changing the stored cell type from the source data.
# (Union "A1,B2" ranges only affect the first area, so set each cell alone.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.287.35"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.932.86"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "251.34"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "0.7231"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("D7").Value = "0.9983"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.3266"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "27.68"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.07166"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D11").Value = "0.7996"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "0.08082"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "1.931.82"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "5.432"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "94.99"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "14.90"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "30.289.15"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "258.02"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "0.000008154"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "5.802"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "2.183.72"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "0.9981"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "0.9990"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "6.882"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "9.654"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "165.41"
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("D27").Value = "19.35"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "2.305"
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").Value = "0.1291"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "1.541"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "4.449"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "4.201"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "0.05237"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").Value = "1.275"
$ws.Range("E35").Value = "  +4.25%  "
$ws.Range("D36").Value = "0.7537"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "2.782"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").Value = "0.01974"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "2.821"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "79.08"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").Value = "6.417"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").Value = "0.4561"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").Value = "2.035"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("D44").Value = "0.8412"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "0.9987"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "101.26"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").Value = "7.486"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").Value = "36.84"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "0.4241"
$ws.Range("E50").Value = "  +3.92%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06085"
$ws.Range("E51").Value = "  +2.12%  "
